$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A from 94 to 96
# (the COM width model adds ~0.83 padding, so use a slightly lower
# input value that still rounds/stores to the target width of 96)
$ws.Columns.Item(1).ColumnWidth = 95.14

# New rows 263-266
$ws.Range("A263").Value = "I want to add this name as my curve shade name length \`"Hydrocarbon bearing zone highlighted\`""
$ws.Range("B263").Value = "llama3.2:latest"
$ws.Range("C263").Value = "To set the curve shading name, click on the `"Curve`" menu and select `"Shading`". Then, in the `"Shading`" dialog box, enter `"Hydrocarbon bearing zone highlighted`" in the `"Name`" field. Note that the character limit for the curve shade name is 20 characters. Since your chosen name has a length of 37 characters, it exceeds the allowed limit. Therefore, this operation is not allowed due to the character length constraint."

$ws.Range("A264").Value = "Why can't I add 251 curve shades to my log?"
$ws.Range("B264").Value = "llama3.2:latest"
$ws.Range("C264").Value = "According to Document 25, the maximum number of curve shades per plot is 250. This means you cannot add more than 250 curve shades to your log."

$ws.Range("A265").Value = "I have 20000 modifiers added ty log, why I can't I add anymore?"
$ws.Range("B265").Value = "llama3.2:latest"
$ws.Range("C265").Value = "You cannot add more than 20000 modifiers per plot because of the limit specified on theHometab."

$ws.Range("A266").Value = "What's the maximum number of data points allowed per curve?"
$ws.Range("B266").Value = "llama3.2:latest"
$ws.Range("C266").Value = "Based on Document 29, which states: `"Deviation from the above may cause errors during loading data from files.`" and considering that the curves are being loaded into columns in Geo so the maximum number of curves correspond to the maximum number of columns in the data file, the answer is:
You can load an unlimited number of data points per curve."

# The multi-line text in C266 makes the engine compute a custom row
# height; AutoFit() puts the row back to the sheet's normal/default
# height so no ht/customHeight attribute is written out.
$ws.Rows.Item(266).AutoFit()
